$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A367").Value = "{'shrinking': True, 'kernel': 'linear', 'gamma': 'auto', 'degree': 3, 'C': 2.03}"
$ws.Range("B367").Value = 0.867
$ws.Range("C367").Value = 0.767

$ws.Range("A368").Value = "{'shrinking': True, 'kernel': 'poly', 'gamma': 'scale', 'degree': 3, 'coef0': 7, 'C': 0.1}"
$ws.Range("B368").Value = 0.858
$ws.Range("C368").Value = 0.767

$ws.Range("A369").Value = "{'C': 0.1, 'coef0': 7, 'degree': 3, 'gamma': 'scale', 'kernel': 'poly', 'shrinking': True}"
$ws.Range("B369").Value = 0.858
$ws.Range("C369").Value = 0.767

$ws.Range("A370").Value = "{'C': 0.1, 'coef0': 6.99999999999994, 'degree': 3, 'gamma': 'scale', 'kernel': 'poly', 'shrinking': True}"
$ws.Range("B370").Value = 0.858
$ws.Range("C370").Value = 0.767

$ws.Range("A371").Value = "{'C': 0.1, 'coef0': 7.0, 'degree': 3, 'gamma': 'scale', 'kernel': 'poly', 'shrinking': True}"
$ws.Range("B371").Value = 0.858
$ws.Range("C371").Value = 0.767

$ws.Range("A372").Value = "{'C': 0.1, 'coef0': 7.0, 'degree': 3, 'gamma': 'scale', 'kernel': 'poly', 'shrinking': True}"
$ws.Range("B372").Value = 0.858
$ws.Range("C372").Value = 0.767
